$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A4: new FP_Card_No regression-test value, pasted in with its own formatting ---
$a4 = $ws.Range("A4")
$a4.ClearFormats()
$a4.NumberFormat = "@"
$a4.Value = "6375004107031383"
$a4.HorizontalAlignment = -4131
$a4.VerticalAlignment = -4160
$a4.WrapText = $true
$a4.Interior.Color = 16777215
$a4.Borders.LineStyle = 1

# --- B4: new FP_Pin value, same pasted-in formatting ---
$b4 = $ws.Range("B4")
$b4.ClearFormats()
$b4.Value = 225
$b4.NumberFormat = "@"
$b4.HorizontalAlignment = -4108
$b4.VerticalAlignment = -4160
$b4.WrapText = $true
$b4.Interior.Color = 16777215
$b4.Borders.LineStyle = 1

$ws.Range("D9").Select()
